# Applies the "Saldo" worksheet update described by the commit diff:
#   - TAYLA (005890232) moves from near the bottom of the list (old
#     balance 79.16, row 118) up to become the new row 3, with an
#     updated balance of 26377.15.
#   - CINTIA's (004927044) balance changes from 19040 to 18871.81.
#   - KELMA's (004504449) row is removed entirely.
#   - BLUEMETRIX (001761119) moves from row 73 (old balance 115.62) to
#     just above THIAGO (row 50), with an updated balance of 283.81.
#
# The sheet has no gaps/merges and rows are addressed purely by position,
# so every row index below is taken from the *original* layout and the
# operations are applied from the bottom of the sheet upward - that way
# each row number used here is still valid at the moment it's used
# (nothing above it has shifted yet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export")

# --- Remove the old TAYLA row (row 118: 005890232 / TAYLA / 79.16) ---
$ws.Rows.Item(118).Delete()

# --- Remove the old BLUEMETRIX row (row 73: 001761119 / BLUEMETRIX / 115.62) ---
$ws.Rows.Item(73).Delete()

# --- Insert the new BLUEMETRIX row right before THIAGO (row 50) ---
$ws.Rows.Item(50).Insert()
$ws.Cells.Item(50, 1).NumberFormat = "@"
$ws.Cells.Item(50, 1).Value = "001761119"
$ws.Cells.Item(50, 2).Value = "BLUEMETRIX"
$ws.Cells.Item(50, 3).Value = 283.81

# --- Remove the KELMA row (row 7: 004504449 / KELMA / 1060.21) ---
$ws.Rows.Item(7).Delete()

# --- Update CINTIA's balance (row 3: 004927044 / CINTIA) ---
$ws.Cells.Item(3, 3).Value = 18871.81

# --- Insert the new TAYLA row right before CINTIA (row 3) ---
$ws.Rows.Item(3).Insert()
$ws.Cells.Item(3, 1).NumberFormat = "@"
$ws.Cells.Item(3, 1).Value = "005890232"
$ws.Cells.Item(3, 2).Value = "TAYLA"
$ws.Cells.Item(3, 3).Value = 26377.15
